$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 560
$ws.Cells.Item(560, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(560, 1).Value = 44291
$ws.Cells.Item(560, 2).NumberFormat = "@"
$ws.Cells.Item(560, 2).Value = '1051155'
$ws.Cells.Item(560, 2).Style = "Normal"
$ws.Cells.Item(560, 3).Value = 3011
$ws.Cells.Item(560, 4).Value = 'Order 1051155 Card(Stripe)'
$ws.Cells.Item(560, 5).Value = "'"
$ws.Cells.Item(560, 5).Style = "Normal"
$ws.Cells.Item(560, 6).Value = 330.36

# Row 561
$ws.Cells.Item(561, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(561, 1).Value = 44291
$ws.Cells.Item(561, 2).NumberFormat = "@"
$ws.Cells.Item(561, 2).Value = '1051155'
$ws.Cells.Item(561, 2).Style = "Normal"
$ws.Cells.Item(561, 3).Value = 2611
$ws.Cells.Item(561, 4).Value = 'Order 1051155 Card(Stripe)'
$ws.Cells.Item(561, 5).Value = "'"
$ws.Cells.Item(561, 5).Style = "Normal"
$ws.Cells.Item(561, 6).Value = 39.64

# Row 562
$ws.Cells.Item(562, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(562, 1).Value = 44291
$ws.Cells.Item(562, 2).NumberFormat = "@"
$ws.Cells.Item(562, 2).Value = '1051155'
$ws.Cells.Item(562, 2).Style = "Normal"
$ws.Cells.Item(562, 3).Value = 1930
$ws.Cells.Item(562, 4).Value = 'Order 1051155 Card(Stripe)'
$ws.Cells.Item(562, 5).Value = 370
$ws.Cells.Item(562, 6).Value = "'"
$ws.Cells.Item(562, 6).Style = "Normal"

# Row 563
$ws.Cells.Item(563, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(563, 1).Value = 44291
$ws.Cells.Item(563, 2).Value = "'"
$ws.Cells.Item(563, 2).Style = "Normal"
$ws.Cells.Item(563, 3).Value = 5670
$ws.Cells.Item(563, 4).Value = 'ST1 V#LLINGBY K6885'
$ws.Cells.Item(563, 5).Value = 663.73
$ws.Cells.Item(563, 6).Value = "'"
$ws.Cells.Item(563, 6).Style = "Normal"

# Row 564
$ws.Cells.Item(564, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(564, 1).Value = 44291
$ws.Cells.Item(564, 2).Value = "'"
$ws.Cells.Item(564, 2).Style = "Normal"
$ws.Cells.Item(564, 3).Value = 2641
$ws.Cells.Item(564, 4).Value = 'ST1 V#LLINGBY K6885'
$ws.Cells.Item(564, 5).Value = 165.93
$ws.Cells.Item(564, 6).Value = "'"
$ws.Cells.Item(564, 6).Style = "Normal"

# Row 565
$ws.Cells.Item(565, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(565, 1).Value = 44291
$ws.Cells.Item(565, 2).Value = "'"
$ws.Cells.Item(565, 2).Style = "Normal"
$ws.Cells.Item(565, 3).Value = 1930
$ws.Cells.Item(565, 4).Value = 'ST1 V#LLINGBY K6885'
$ws.Cells.Item(565, 5).Value = "'"
$ws.Cells.Item(565, 5).Style = "Normal"
$ws.Cells.Item(565, 6).Value = 829.66

# Row 566
$ws.Cells.Item(566, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(566, 1).Value = 44292
$ws.Cells.Item(566, 2).Value = 'Reko75'
$ws.Cells.Item(566, 3).Value = 3011
$ws.Cells.Item(566, 4).Value = 'Reko Swish +46733035539'
$ws.Cells.Item(566, 5).Value = "'"
$ws.Cells.Item(566, 5).Style = "Normal"
$ws.Cells.Item(566, 6).Value = 460.71

# Row 567
$ws.Cells.Item(567, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(567, 1).Value = 44292
$ws.Cells.Item(567, 2).Value = 'Reko75'
$ws.Cells.Item(567, 3).Value = 2611
$ws.Cells.Item(567, 4).Value = 'Reko Swish +46733035539'
$ws.Cells.Item(567, 5).Value = "'"
$ws.Cells.Item(567, 5).Style = "Normal"
$ws.Cells.Item(567, 6).Value = 55.29

# Row 568
$ws.Cells.Item(568, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(568, 1).Value = 44292
$ws.Cells.Item(568, 2).Value = 'Reko75'
$ws.Cells.Item(568, 3).Value = 1930
$ws.Cells.Item(568, 4).Value = 'Reko Swish +46733035539'
$ws.Cells.Item(568, 5).Value = 516
$ws.Cells.Item(568, 6).Value = "'"
$ws.Cells.Item(568, 6).Style = "Normal"

# Row 569
$ws.Cells.Item(569, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(569, 1).Value = 44292
$ws.Cells.Item(569, 2).Value = 'Reko76'
$ws.Cells.Item(569, 3).Value = 3011
$ws.Cells.Item(569, 4).Value = 'Reko Swish +46731835553'
$ws.Cells.Item(569, 5).Value = "'"
$ws.Cells.Item(569, 5).Style = "Normal"
$ws.Cells.Item(569, 6).Value = 185.71

# Row 570
$ws.Cells.Item(570, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(570, 1).Value = 44292
$ws.Cells.Item(570, 2).Value = 'Reko76'
$ws.Cells.Item(570, 3).Value = 2611
$ws.Cells.Item(570, 4).Value = 'Reko Swish +46731835553'
$ws.Cells.Item(570, 5).Value = "'"
$ws.Cells.Item(570, 5).Style = "Normal"
$ws.Cells.Item(570, 6).Value = 22.29

# Row 571
$ws.Cells.Item(571, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(571, 1).Value = 44292
$ws.Cells.Item(571, 2).Value = 'Reko76'
$ws.Cells.Item(571, 3).Value = 1930
$ws.Cells.Item(571, 4).Value = 'Reko Swish +46731835553'
$ws.Cells.Item(571, 5).Value = 208
$ws.Cells.Item(571, 6).Value = "'"
$ws.Cells.Item(571, 6).Style = "Normal"

# Row 572
$ws.Cells.Item(572, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(572, 1).Value = 44292
$ws.Cells.Item(572, 2).NumberFormat = "@"
$ws.Cells.Item(572, 2).Value = '9061234'
$ws.Cells.Item(572, 2).Style = "Normal"
$ws.Cells.Item(572, 3).Value = 3011
$ws.Cells.Item(572, 4).Value = 'Order 9061234 Swish +46705293845'
$ws.Cells.Item(572, 5).Value = "'"
$ws.Cells.Item(572, 5).Style = "Normal"
$ws.Cells.Item(572, 6).Value = 610.71

# Row 573
$ws.Cells.Item(573, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(573, 1).Value = 44292
$ws.Cells.Item(573, 2).NumberFormat = "@"
$ws.Cells.Item(573, 2).Value = '9061234'
$ws.Cells.Item(573, 2).Style = "Normal"
$ws.Cells.Item(573, 3).Value = 2611
$ws.Cells.Item(573, 4).Value = 'Order 9061234 Swish +46705293845'
$ws.Cells.Item(573, 5).Value = "'"
$ws.Cells.Item(573, 5).Style = "Normal"
$ws.Cells.Item(573, 6).Value = 73.29

# Row 574
$ws.Cells.Item(574, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(574, 1).Value = 44292
$ws.Cells.Item(574, 2).NumberFormat = "@"
$ws.Cells.Item(574, 2).Value = '9061234'
$ws.Cells.Item(574, 2).Style = "Normal"
$ws.Cells.Item(574, 3).Value = 1930
$ws.Cells.Item(574, 4).Value = 'Order 9061234 Swish +46705293845'
$ws.Cells.Item(574, 5).Value = 684
$ws.Cells.Item(574, 6).Value = "'"
$ws.Cells.Item(574, 6).Style = "Normal"

# Row 575
$ws.Cells.Item(575, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(575, 1).Value = 44292
$ws.Cells.Item(575, 2).Value = 'Reko77'
$ws.Cells.Item(575, 3).Value = 3011
$ws.Cells.Item(575, 4).Value = 'Reko Swish +46709622907'
$ws.Cells.Item(575, 5).Value = "'"
$ws.Cells.Item(575, 5).Style = "Normal"
$ws.Cells.Item(575, 6).Value = 282.14

# Row 576
$ws.Cells.Item(576, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(576, 1).Value = 44292
$ws.Cells.Item(576, 2).Value = 'Reko77'
$ws.Cells.Item(576, 3).Value = 2611
$ws.Cells.Item(576, 4).Value = 'Reko Swish +46709622907'
$ws.Cells.Item(576, 5).Value = "'"
$ws.Cells.Item(576, 5).Style = "Normal"
$ws.Cells.Item(576, 6).Value = 33.86

# Row 577
$ws.Cells.Item(577, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(577, 1).Value = 44292
$ws.Cells.Item(577, 2).Value = 'Reko77'
$ws.Cells.Item(577, 3).Value = 1930
$ws.Cells.Item(577, 4).Value = 'Reko Swish +46709622907'
$ws.Cells.Item(577, 5).Value = 316
$ws.Cells.Item(577, 6).Value = "'"
$ws.Cells.Item(577, 6).Style = "Normal"

# Row 578
$ws.Cells.Item(578, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(578, 1).Value = 44293
$ws.Cells.Item(578, 2).Value = 'Reko78'
$ws.Cells.Item(578, 3).Value = 3011
$ws.Cells.Item(578, 4).Value = 'Reko Swish +46709906521'
$ws.Cells.Item(578, 5).Value = "'"
$ws.Cells.Item(578, 5).Style = "Normal"
$ws.Cells.Item(578, 6).Value = 141.07

# Row 579
$ws.Cells.Item(579, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(579, 1).Value = 44293
$ws.Cells.Item(579, 2).Value = 'Reko78'
$ws.Cells.Item(579, 3).Value = 2611
$ws.Cells.Item(579, 4).Value = 'Reko Swish +46709906521'
$ws.Cells.Item(579, 5).Value = "'"
$ws.Cells.Item(579, 5).Style = "Normal"
$ws.Cells.Item(579, 6).Value = 16.93

# Row 580
$ws.Cells.Item(580, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(580, 1).Value = 44293
$ws.Cells.Item(580, 2).Value = 'Reko78'
$ws.Cells.Item(580, 3).Value = 1930
$ws.Cells.Item(580, 4).Value = 'Reko Swish +46709906521'
$ws.Cells.Item(580, 5).Value = 158
$ws.Cells.Item(580, 6).Value = "'"
$ws.Cells.Item(580, 6).Style = "Normal"

# Row 581
$ws.Cells.Item(581, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(581, 1).Value = 44293
$ws.Cells.Item(581, 2).NumberFormat = "@"
$ws.Cells.Item(581, 2).Value = '5072119'
$ws.Cells.Item(581, 2).Style = "Normal"
$ws.Cells.Item(581, 3).Value = 3011
$ws.Cells.Item(581, 4).Value = 'Order 5072119 Swish +46730402047'
$ws.Cells.Item(581, 5).Value = "'"
$ws.Cells.Item(581, 5).Style = "Normal"
$ws.Cells.Item(581, 6).Value = 352.68

# Row 582
$ws.Cells.Item(582, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(582, 1).Value = 44293
$ws.Cells.Item(582, 2).NumberFormat = "@"
$ws.Cells.Item(582, 2).Value = '5072119'
$ws.Cells.Item(582, 2).Style = "Normal"
$ws.Cells.Item(582, 3).Value = 2611
$ws.Cells.Item(582, 4).Value = 'Order 5072119 Swish +46730402047'
$ws.Cells.Item(582, 5).Value = "'"
$ws.Cells.Item(582, 5).Style = "Normal"
$ws.Cells.Item(582, 6).Value = 42.32

# Row 583
$ws.Cells.Item(583, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(583, 1).Value = 44293
$ws.Cells.Item(583, 2).NumberFormat = "@"
$ws.Cells.Item(583, 2).Value = '5072119'
$ws.Cells.Item(583, 2).Style = "Normal"
$ws.Cells.Item(583, 3).Value = 1930
$ws.Cells.Item(583, 4).Value = 'Order 5072119 Swish +46730402047'
$ws.Cells.Item(583, 5).Value = 395
$ws.Cells.Item(583, 6).Value = "'"
$ws.Cells.Item(583, 6).Style = "Normal"

# Row 584
$ws.Cells.Item(584, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(584, 1).Value = 44294
$ws.Cells.Item(584, 2).NumberFormat = "@"
$ws.Cells.Item(584, 2).Value = '3081627'
$ws.Cells.Item(584, 2).Style = "Normal"
$ws.Cells.Item(584, 3).Value = 3011
$ws.Cells.Item(584, 4).Value = 'Order 3081627 Swish +46768551925'
$ws.Cells.Item(584, 5).Value = "'"
$ws.Cells.Item(584, 5).Style = "Normal"
$ws.Cells.Item(584, 6).Value = 956.25

# Row 585
$ws.Cells.Item(585, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(585, 1).Value = 44294
$ws.Cells.Item(585, 2).NumberFormat = "@"
$ws.Cells.Item(585, 2).Value = '3081627'
$ws.Cells.Item(585, 2).Style = "Normal"
$ws.Cells.Item(585, 3).Value = 2611
$ws.Cells.Item(585, 4).Value = 'Order 3081627 Swish +46768551925'
$ws.Cells.Item(585, 5).Value = "'"
$ws.Cells.Item(585, 5).Style = "Normal"
$ws.Cells.Item(585, 6).Value = 114.75

# Row 586
$ws.Cells.Item(586, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(586, 1).Value = 44294
$ws.Cells.Item(586, 2).NumberFormat = "@"
$ws.Cells.Item(586, 2).Value = '3081627'
$ws.Cells.Item(586, 2).Style = "Normal"
$ws.Cells.Item(586, 3).Value = 1930
$ws.Cells.Item(586, 4).Value = 'Order 3081627 Swish +46768551925'
$ws.Cells.Item(586, 5).Value = 1071
$ws.Cells.Item(586, 6).Value = "'"
$ws.Cells.Item(586, 6).Style = "Normal"

# Row 587
$ws.Cells.Item(587, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(587, 1).Value = 44294
$ws.Cells.Item(587, 2).NumberFormat = "@"
$ws.Cells.Item(587, 2).Value = '7082223'
$ws.Cells.Item(587, 2).Style = "Normal"
$ws.Cells.Item(587, 3).Value = 3011
$ws.Cells.Item(587, 4).Value = 'Order 7082223 Swish +46723656673'
$ws.Cells.Item(587, 5).Value = "'"
$ws.Cells.Item(587, 5).Style = "Normal"
$ws.Cells.Item(587, 6).Value = 806.25

# Row 588
$ws.Cells.Item(588, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(588, 1).Value = 44294
$ws.Cells.Item(588, 2).NumberFormat = "@"
$ws.Cells.Item(588, 2).Value = '7082223'
$ws.Cells.Item(588, 2).Style = "Normal"
$ws.Cells.Item(588, 3).Value = 2611
$ws.Cells.Item(588, 4).Value = 'Order 7082223 Swish +46723656673'
$ws.Cells.Item(588, 5).Value = "'"
$ws.Cells.Item(588, 5).Style = "Normal"
$ws.Cells.Item(588, 6).Value = 96.75

# Row 589
$ws.Cells.Item(589, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(589, 1).Value = 44294
$ws.Cells.Item(589, 2).NumberFormat = "@"
$ws.Cells.Item(589, 2).Value = '7082223'
$ws.Cells.Item(589, 2).Style = "Normal"
$ws.Cells.Item(589, 3).Value = 1930
$ws.Cells.Item(589, 4).Value = 'Order 7082223 Swish +46723656673'
$ws.Cells.Item(589, 5).Value = 903
$ws.Cells.Item(589, 6).Value = "'"
$ws.Cells.Item(589, 6).Style = "Normal"

# Row 590
$ws.Cells.Item(590, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(590, 1).Value = 44294
$ws.Cells.Item(590, 2).Value = "'"
$ws.Cells.Item(590, 2).Style = "Normal"
$ws.Cells.Item(590, 3).Value = 6400
$ws.Cells.Item(590, 4).Value = 'FACEBK FZDBG3KZ62 K6885'
$ws.Cells.Item(590, 5).Value = 430
$ws.Cells.Item(590, 6).Value = "'"
$ws.Cells.Item(590, 6).Style = "Normal"

# Row 591
$ws.Cells.Item(591, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(591, 1).Value = 44294
$ws.Cells.Item(591, 2).Value = "'"
$ws.Cells.Item(591, 2).Style = "Normal"
$ws.Cells.Item(591, 3).Value = "'"
$ws.Cells.Item(591, 3).Style = "Normal"
$ws.Cells.Item(591, 4).Value = 'FACEBK FZDBG3KZ62 K6885'
$ws.Cells.Item(591, 5).Value = 0
$ws.Cells.Item(591, 6).Value = "'"
$ws.Cells.Item(591, 6).Style = "Normal"

# Row 592
$ws.Cells.Item(592, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(592, 1).Value = 44294
$ws.Cells.Item(592, 2).Value = "'"
$ws.Cells.Item(592, 2).Style = "Normal"
$ws.Cells.Item(592, 3).Value = 1930
$ws.Cells.Item(592, 4).Value = 'FACEBK FZDBG3KZ62 K6885'
$ws.Cells.Item(592, 5).Value = "'"
$ws.Cells.Item(592, 5).Style = "Normal"
$ws.Cells.Item(592, 6).Value = 430

# Row 593
$ws.Cells.Item(593, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(593, 1).Value = 44295
$ws.Cells.Item(593, 2).Value = 'Reko79'
$ws.Cells.Item(593, 3).Value = 3011
$ws.Cells.Item(593, 4).Value = 'Reko Swish +46703344337'
$ws.Cells.Item(593, 5).Value = "'"
$ws.Cells.Item(593, 5).Style = "Normal"
$ws.Cells.Item(593, 6).Value = 345.54

# Row 594
$ws.Cells.Item(594, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(594, 1).Value = 44295
$ws.Cells.Item(594, 2).Value = 'Reko79'
$ws.Cells.Item(594, 3).Value = 2611
$ws.Cells.Item(594, 4).Value = 'Reko Swish +46703344337'
$ws.Cells.Item(594, 5).Value = "'"
$ws.Cells.Item(594, 5).Style = "Normal"
$ws.Cells.Item(594, 6).Value = 41.46

# Row 595
$ws.Cells.Item(595, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(595, 1).Value = 44295
$ws.Cells.Item(595, 2).Value = 'Reko79'
$ws.Cells.Item(595, 3).Value = 1930
$ws.Cells.Item(595, 4).Value = 'Reko Swish +46703344337'
$ws.Cells.Item(595, 5).Value = 387
$ws.Cells.Item(595, 6).Value = "'"
$ws.Cells.Item(595, 6).Style = "Normal"

# Row 596
$ws.Cells.Item(596, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(596, 1).Value = 44295
$ws.Cells.Item(596, 2).Value = "'"
$ws.Cells.Item(596, 2).Style = "Normal"
$ws.Cells.Item(596, 3).Value = 4010
$ws.Cells.Item(596, 4).Value = 'M&S RB BROMMA K0135'
$ws.Cells.Item(596, 5).Value = 1929.28
$ws.Cells.Item(596, 6).Value = "'"
$ws.Cells.Item(596, 6).Style = "Normal"

# Row 597
$ws.Cells.Item(597, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(597, 1).Value = 44295
$ws.Cells.Item(597, 2).Value = "'"
$ws.Cells.Item(597, 2).Style = "Normal"
$ws.Cells.Item(597, 3).Value = 2645
$ws.Cells.Item(597, 4).Value = 'M&S RB BROMMA K0135'
$ws.Cells.Item(597, 5).Value = 231.51
$ws.Cells.Item(597, 6).Value = "'"
$ws.Cells.Item(597, 6).Style = "Normal"

# Row 598
$ws.Cells.Item(598, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(598, 1).Value = 44295
$ws.Cells.Item(598, 2).Value = "'"
$ws.Cells.Item(598, 2).Style = "Normal"
$ws.Cells.Item(598, 3).Value = 1930
$ws.Cells.Item(598, 4).Value = 'M&S RB BROMMA K0135'
$ws.Cells.Item(598, 5).Value = "'"
$ws.Cells.Item(598, 5).Style = "Normal"
$ws.Cells.Item(598, 6).Value = 2160.79

# Row 599
$ws.Cells.Item(599, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(599, 1).Value = 44296
$ws.Cells.Item(599, 2).NumberFormat = "@"
$ws.Cells.Item(599, 2).Value = '6101315'
$ws.Cells.Item(599, 2).Style = "Normal"
$ws.Cells.Item(599, 3).Value = 3011
$ws.Cells.Item(599, 4).Value = 'Order 6101315 Card(Stripe)'
$ws.Cells.Item(599, 5).Value = "'"
$ws.Cells.Item(599, 5).Style = "Normal"
$ws.Cells.Item(599, 6).Value = 1062.5

# Row 600
$ws.Cells.Item(600, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(600, 1).Value = 44296
$ws.Cells.Item(600, 2).NumberFormat = "@"
$ws.Cells.Item(600, 2).Value = '6101315'
$ws.Cells.Item(600, 2).Style = "Normal"
$ws.Cells.Item(600, 3).Value = 2611
$ws.Cells.Item(600, 4).Value = 'Order 6101315 Card(Stripe)'
$ws.Cells.Item(600, 5).Value = "'"
$ws.Cells.Item(600, 5).Style = "Normal"
$ws.Cells.Item(600, 6).Value = 127.5

# Row 601
$ws.Cells.Item(601, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(601, 1).Value = 44296
$ws.Cells.Item(601, 2).NumberFormat = "@"
$ws.Cells.Item(601, 2).Value = '6101315'
$ws.Cells.Item(601, 2).Style = "Normal"
$ws.Cells.Item(601, 3).Value = 1930
$ws.Cells.Item(601, 4).Value = 'Order 6101315 Card(Stripe)'
$ws.Cells.Item(601, 5).Value = 1190
$ws.Cells.Item(601, 6).Value = "'"
$ws.Cells.Item(601, 6).Style = "Normal"

# Row 602
$ws.Cells.Item(602, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(602, 1).Value = 44297
$ws.Cells.Item(602, 2).Value = "'"
$ws.Cells.Item(602, 2).Style = "Normal"
$ws.Cells.Item(602, 3).Value = 4010
$ws.Cells.Item(602, 4).Value = 'NGROCERIES K0135'
$ws.Cells.Item(602, 5).Value = 176.79
$ws.Cells.Item(602, 6).Value = "'"
$ws.Cells.Item(602, 6).Style = "Normal"

# Row 603
$ws.Cells.Item(603, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(603, 1).Value = 44297
$ws.Cells.Item(603, 2).Value = "'"
$ws.Cells.Item(603, 2).Style = "Normal"
$ws.Cells.Item(603, 3).Value = 2645
$ws.Cells.Item(603, 4).Value = 'NGROCERIES K0135'
$ws.Cells.Item(603, 5).Value = 21.21
$ws.Cells.Item(603, 6).Value = "'"
$ws.Cells.Item(603, 6).Style = "Normal"

# Row 604
$ws.Cells.Item(604, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(604, 1).Value = 44297
$ws.Cells.Item(604, 2).Value = "'"
$ws.Cells.Item(604, 2).Style = "Normal"
$ws.Cells.Item(604, 3).Value = 1930
$ws.Cells.Item(604, 4).Value = 'NGROCERIES K0135'
$ws.Cells.Item(604, 5).Value = "'"
$ws.Cells.Item(604, 5).Style = "Normal"
$ws.Cells.Item(604, 6).Value = 198

